# Update the weekly Fruta/Hortaliza price data for rows 2-13.
# Columns changed per row: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44453; J = 20;  K = 2300; L = 2300; M = 2300; P = 2300 }
    3  = @{ D = 44497; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    4  = @{ D = 44496; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
    5  = @{ D = 44484; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
    6  = @{ D = 44203; J = 30;  K = 2000; L = 2000; M = 2000; P = 2000 }
    7  = @{ D = 44487; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    8  = @{ D = 44483; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    9  = @{ D = 44476; J = 30;  K = 2200; L = 2200; M = 2200; P = 2200 }
    10 = @{ D = 44447; J = 75;  K = 2200; L = 2200; M = 2200; P = 2200 }
    11 = @{ D = 44474; J = 20;  K = 1600; L = 1600; M = 1600; P = 1600 }
    12 = @{ D = 44452; J = 120; K = 2300; L = 2300; M = 2300; P = 2300 }
    13 = @{ D = 44473; J = 140; K = 1600; L = 1600; M = 1600; P = 1600 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value2  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value2 = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value2 = $vals.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $vals.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value2 = $vals.P   # P: Precio $/Kg
}
